$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New problem block: "count no of subarr xor k"
$ws.Range("A76").Value2 = "count no of subarr xor k"
$ws.Range("C76").Value2 = "@ use 2 for loop to find all possible subarr and find xor val"
$ws.Range("D76").Value2 = "@ use prefix xor and store the values and count of appearance as key"
$ws.Range("E76").Value2 = "@ the target subarr in arr"

$ws.Range("C77").Value2 = "@ takes O(n2)time"
$ws.Range("D77").Value2 = "@ now if prefix xored to target gives me a value which is present in dict then it means there contains"
$ws.Range("E77").Value2 = "@ the count is used to get how many times it appeared"

# Update the view to scroll down to the newly added rows and move the
# active selection to C78, matching the author's final cursor position.
$win = $excel.ActiveWindow
$win.ScrollRow = 65
$win.ScrollColumn = 3
[void]$ws.Range("C78").Select()
